$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 9 (the extra data rows no longer needed)
$ws.Rows("3:9").Delete()

# Update the single remaining data row: B2 keeps location "A21" (was merged from
# old row 8), add Producto/Marca values in the new columns
$ws.Range("B2").Value = "A21"
$ws.Range("C2").Value = "Diego"
$ws.Range("D2").Value = "PP"

# Add the new headers, matching the look of the existing header cells
$ws.Range("C1").Value = "Producto"
$ws.Range("D1").Value = "Marca"

$headerSample = $ws.Range("A1")
$newHeaders = $ws.Range("C1:D1")
$newHeaders.Font.Bold = $headerSample.Font.Bold
$newHeaders.HorizontalAlignment = $headerSample.HorizontalAlignment
$newHeaders.VerticalAlignment = $headerSample.VerticalAlignment
$newHeaders.Borders.LineStyle = $headerSample.Borders.LineStyle
